$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '37.720.95'
$ws.Cells.Item(2, 5).Value = '  +0.62%  '
$ws.Cells.Item(3, 4).Value = '2.023.97'
$ws.Cells.Item(3, 5).Value = '  -0.47%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '226.99'
$ws.Cells.Item(5, 5).Value = '  -1.13%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.605'
$ws.Cells.Item(6, 5).Value = '  -1.24%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '59.82'
$ws.Cells.Item(7, 5).Value = '  +6.73%  '
$ws.Cells.Item(8, 5).Value = '  -0.01%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.380'
$ws.Cells.Item(9, 5).Value = '  -0.41%  '
$ws.Cells.Item(10, 5).Value = '  +1.04%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.104'
$ws.Cells.Item(11, 5).Value = '  +0.44%  '
$ws.Cells.Item(12, 2).Value = 'Chainlink'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '14.54'
$ws.Cells.Item(12, 5).Value = '  +1.00%  '
$ws.Cells.Item(13, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(13, 4).Value = '2.322.21'
$ws.Cells.Item(13, 5).Value = '  -0.50%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '20.96'
$ws.Cells.Item(14, 5).Value = '  +3.11%  '
$ws.Cells.Item(15, 5).Value = '  +1.02%  '
$ws.Cells.Item(16, 5).Value = '  -0.04%  '
$ws.Cells.Item(17, 4).Value = '2.032.24'
$ws.Cells.Item(17, 5).Value = '  +0.45%  '
$ws.Cells.Item(18, 4).Value = '37.667.90'
$ws.Cells.Item(18, 5).Value = '  +0.66%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.10'
$ws.Cells.Item(19, 5).Value = '  -1.68%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '69.43'
$ws.Cells.Item(20, 5).Value = '  +0.57%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0823'
$ws.Cells.Item(21, 5).Value = '  -0.01%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '223.50'
$ws.Cells.Item(22, 5).Value = '  -0.04%  '
$ws.Cells.Item(23, 5).Value = '  +0.05%  '
$ws.Cells.Item(24, 5).Value = '  -1.07%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.19'
$ws.Cells.Item(25, 5).Value = '  -2.91%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '165.72'
$ws.Cells.Item(26, 5).Value = '  +0.43%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9.10'
$ws.Cells.Item(27, 5).Value = '  -0.39%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.127'
$ws.Cells.Item(28, 5).Value = '  -3.39%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '18.79'
$ws.Cells.Item(29, 5).Value = '  +0.11%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.28'
$ws.Cells.Item(30, 5).Value = '  -3.28%  '
$ws.Cells.Item(31, 5).Value = '  +1.37%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.43'
$ws.Cells.Item(32, 5).Value = '  -1.27%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '2.04'
$ws.Cells.Item(33, 5).Value = '  +2.08%  '
$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0602'
$ws.Cells.Item(34, 5).Value = '  -0.85%  '
$ws.Cells.Item(35, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '4.49'
$ws.Cells.Item(35, 5).Value = '  +0.14%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '6.23'
$ws.Cells.Item(36, 5).Value = '  +8.17%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.25'
$ws.Cells.Item(37, 5).Value = '  -2.97%  '
$ws.Cells.Item(38, 5).Value = '  -1.25%  '
$ws.Cells.Item(39, 5).Value = '  -0.24%  '
$ws.Cells.Item(40, 4).Value = '1.529.84'
$ws.Cells.Item(40, 5).Value = '  +3.71%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0216'
$ws.Cells.Item(41, 5).Value = '  +0.37%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '96.27'
$ws.Cells.Item(42, 5).Value = '  +1.36%  '
$ws.Cells.Item(43, 2).Value = 'HuobiToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.81'
$ws.Cells.Item(43, 5).Value = '  -0.74%  '
$ws.Cells.Item(44, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '16.42'
$ws.Cells.Item(44, 5).Value = '  +0.38%  '
$ws.Cells.Item(45, 5).Value = '  -1.52%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.10'
$ws.Cells.Item(46, 5).Value = '  -0.74%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.99'
$ws.Cells.Item(47, 5).Value = '  -5.32%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.95'
$ws.Cells.Item(48, 5).Value = '  +0.29%  '
$ws.Cells.Item(49, 5).Value = '  -1.10%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '7.01'
$ws.Cells.Item(50, 5).Value = '  -1.53%  '
$ws.Cells.Item(51, 4).Value = '2.213.76'
$ws.Cells.Item(51, 5).Value = '  -0.36%  '
